# Fix some failing tests.
#
# RecordInstanceHorizontal.xlsx: D2 held a date serial (44632 / 2022-03-12)
# formatted with a date number format. The fix stores the date as the
# literal text "2022-03-12" instead, using the sheet's default (General)
# style - matching header/body text like the rest of the row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")

# Force the incoming literal to be interpreted as text rather than being
# auto-recognised as a date (which is what a plain .Value assignment of
# "2022-03-12" would otherwise do).
$cell.NumberFormat = "@"
$cell.Value = "2022-03-12"

# Drop the cell back to the default/general format so it matches the
# plain text styling used elsewhere in the row (no leftover custom
# number format on the cell).
$cell.NumberFormat = ""
